$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Market Data")
$ws2 = $wb.Worksheets.Item("Analysis")

# --- Market Data sheet updates ---
$ws1.Range("C2").Value = 97723.44
$ws1.Range("D2").Value = 1493445441.35651
$ws1.Range("E2").Value = 15282.36666
$ws1.Range("F2").Value = 0.354
$ws1.Range("G2").Value = 97689.78722032
$ws1.Range("C3").Value = 2685.93
$ws1.Range("D3").Value = 997614352.494786
$ws1.Range("E3").Value = 371422.3202
$ws1.Range("F3").Value = 1.582
$ws1.Range("G3").Value = 2684.33109043
$ws1.Range("C4").Value = 0.9997
$ws1.Range("D4").Value = 902809385.9070001
$ws1.Range("E4").Value = 903080310
$ws1.Range("F4").Value = -0.01
$ws1.Range("G4").Value = 0.99982084
$ws1.Range("C5").Value = 201.89
$ws1.Range("D5").Value = 857313077.70407
$ws1.Range("E5").Value = 4246436.563
$ws1.Range("F5").Value = -0.694
$ws1.Range("G5").Value = 203.48269591
$ws1.Range("C6").Value = 0.9991
$ws1.Range("D6").Value = 511622726.7582
$ws1.Range("E6").Value = 512083602
$ws1.Range("F6").Value = 0.01
$ws1.Range("G6").Value = 0.99913401
$ws1.Range("C7").Value = 2.4785
$ws1.Range("D7").Value = 493634861.7515
$ws1.Range("E7").Value = 199166779
$ws1.Range("F7").Value = 2.6
$ws1.Range("G7").Value = 2.46397955
$ws1.Range("C8").Value = 0.2117
$ws1.Range("D8").Value = 352624169.23118
$ws1.Range("E8").Value = 1665678645.4
$ws1.Range("F8").Value = 33.48
$ws1.Range("G8").Value = 0.19521403
$ws1.Range("C9").Value = 129.53
$ws1.Range("D9").Value = 292963475.06733
$ws1.Range("E9").Value = 2261742.261
$ws1.Range("F9").Value = 11.051
$ws1.Range("G9").Value = 124.17140339
$ws1.Range("C10").Value = 3.467
$ws1.Range("D10").Value = 286232556.4084
$ws1.Range("E10").Value = 82559145.2
$ws1.Range("F10").Value = 8.378
$ws1.Range("G10").Value = 3.36503685
$ws1.Range("C11").Value = 638.64
$ws1.Range("D11").Value = 263555275.90752
$ws1.Range("E11").Value = 412682.068
$ws1.Range("F11").Value = 5.625
$ws1.Range("G11").Value = 626.13867103
$ws1.Range("C12").Value = 0.26276
$ws1.Range("D12").Value = 236436550.27368
$ws1.Range("E12").Value = 899819418
$ws1.Range("F12").Value = 5.121
$ws1.Range("G12").Value = 0.25954777
$ws1.Range("C13").Value = 0.6788
$ws1.Range("D13").Value = 225714302.71512
$ws1.Range("E13").Value = 332519597.4
$ws1.Range("F13").Value = 35.3
$ws1.Range("G13").Value = 0.66187316
$ws1.Range("C14").Value = 0.00001023
$ws1.Range("D14").Value = 219520535.6558877
$ws1.Range("E14").Value = 21458507884251
$ws1.Range("F14").Value = 7.458
$ws1.Range("G14").Value = 0.00000999
$ws1.Range("C15").Value = 0.8053
$ws1.Range("D15").Value = 226602279.12023
$ws1.Range("E15").Value = 281388649.1
$ws1.Range("F15").Value = 15.555
$ws1.Range("G15").Value = 0.75739919
$ws1.Range("C16").Value = 15.88
$ws1.Range("D16").Value = 188823230.26852
$ws1.Range("E16").Value = 11890631.629
$ws1.Range("F16").Value = -0.626
$ws1.Range("G16").Value = 16.21612462
$ws1.Range("D17").Value = 125586081.88822
$ws1.Range("E17").Value = 505987437.1
$ws1.Range("F17").Value = 4.023
$ws1.Range("G17").Value = 0.24399162
$ws1.Range("C18").Value = 0.4754
$ws1.Range("D18").Value = 103706108.550178
$ws1.Range("E18").Value = 218144948.57
$ws1.Range("F18").Value = -0.979
$ws1.Range("G18").Value = 0.48801697
$ws1.Range("C19").Value = 0.2847
$ws1.Range("D19").Value = 98308574.18538
$ws1.Range("E19").Value = 345305845.4
$ws1.Range("F19").Value = 11.081
$ws1.Range("G19").Value = 0.28683605
$ws1.Range("C20").Value = 1.299
$ws1.Range("D20").Value = 86860919.3916
$ws1.Range("E20").Value = 66867528.4
$ws1.Range("F20").Value = 2.203
$ws1.Range("G20").Value = 1.30058832
$ws1.Range("C21").Value = 0.001093
$ws1.Range("D21").Value = 90251814.621605
$ws1.Range("E21").Value = 82572565985
$ws1.Range("F21").Value = 61.686
$ws1.Range("G21").Value = 0.00096763
$ws1.Range("A22").Value = 'The Anthropic Order'
$ws1.Range("B22").Value = 'TAO'
$ws1.Range("C22").Value = 421.2
$ws1.Range("D22").Value = 79876419.63912
$ws1.Range("E22").Value = 189640.1226
$ws1.Range("F22").Value = 8.305
$ws1.Range("G22").Value = 408.1307934
$ws1.Range("A23").Value = 'Berachain'
$ws1.Range("B23").Value = 'BERA'
$ws1.Range("C23").Value = 5.57
$ws1.Range("D23").Value = 74940461.82808
$ws1.Range("E23").Value = 13454301.944
$ws1.Range("F23").Value = 11.892
$ws1.Range("G23").Value = 5.58017001
$ws1.Range("A24").Value = 'PancakeSwap'
$ws1.Range("B24").Value = 'CAKE'
$ws1.Range("C24").Value = 1.978
$ws1.Range("D24").Value = 78365790.83514
$ws1.Range("E24").Value = 39618701.13
$ws1.Range("F24").Value = 19.444
$ws1.Range("G24").Value = 1.89254516
$ws1.Range("A25").Value = 'Aptos'
$ws1.Range("B25").Value = 'APT'
$ws1.Range("C25").Value = 6.08
$ws1.Range("D25").Value = 68274945.04
$ws1.Range("E25").Value = 11229431.75
$ws1.Range("F25").Value = -0.977
$ws1.Range("G25").Value = 6.24965599
$ws1.Range("C26").Value = 0.239
$ws1.Range("D26").Value = 66168023.8557
$ws1.Range("E26").Value = 276853656.3
$ws1.Range("F26").Value = 19.381
$ws1.Range("G26").Value = 0.2532168
$ws1.Range("A27").Value = 'THORChain'
$ws1.Range("B27").Value = 'RUNE'
$ws1.Range("C27").Value = 1.401
$ws1.Range("D27").Value = 68714880.0636
$ws1.Range("E27").Value = 49047023.6
$ws1.Range("F27").Value = 1.228
$ws1.Range("G27").Value = 1.41607696
$ws1.Range("C28").Value = 19.38
$ws1.Range("D28").Value = 67234017.519
$ws1.Range("E28").Value = 3469247.55
$ws1.Range("F28").Value = 4.026
$ws1.Range("G28").Value = 19.15883459
$ws1.Range("C29").Value = 0.24018
$ws1.Range("D29").Value = 65921330.3448
$ws1.Range("E29").Value = 274466360
$ws1.Range("F29").Value = 1.44
$ws1.Range("G29").Value = 0.23936172
$ws1.Range("C30").Value = 0.2645
$ws1.Range("D30").Value = 61744974.50585001
$ws1.Range("E30").Value = 233440357.3
$ws1.Range("F30").Value = 49.774
$ws1.Range("G30").Value = 0.26524064
$ws1.Range("C31").Value = 0.657
$ws1.Range("D31").Value = 56550517.72866
$ws1.Range("E31").Value = 86073847.38
$ws1.Range("F31").Value = -0.152
$ws1.Range("G31").Value = 0.67329129
$ws1.Range("A32").Value = 'Aave'
$ws1.Range("B32").Value = 'AAVE'
$ws1.Range("C32").Value = 255.37
$ws1.Range("D32").Value = 54302407.71423
$ws1.Range("E32").Value = 212642.079
$ws1.Range("F32").Value = 3.13
$ws1.Range("G32").Value = 255.94602405
$ws1.Range("A33").Value = 'MANTRA'
$ws1.Range("B33").Value = 'OM'
$ws1.Range("C33").Value = 5.9791
$ws1.Range("D33").Value = 51369909.9989
$ws1.Range("E33").Value = 8591579
$ws1.Range("F33").Value = -2.9
$ws1.Range("G33").Value = 6.0390153
$ws1.Range("C34").Value = 0.1609
$ws1.Range("D34").Value = 48194222.37509999
$ws1.Range("E34").Value = 299529039
$ws1.Range("F34").Value = -2.838
$ws1.Range("G34").Value = 0.17213671
$ws1.Range("A35").Value = 'Meld Bridged AVAX (Meld)'
$ws1.Range("B35").Value = 'AVAX'
$ws1.Range("C35").Value = 26.36
$ws1.Range("D35").Value = 48034086.922
$ws1.Range("E35").Value = 1822233.95
$ws1.Range("F35").Value = 4.272
$ws1.Range("G35").Value = 26.04792904
$ws1.Range("A36").Value = 'NEAR Protocol'
$ws1.Range("B36").Value = 'NEAR'
$ws1.Range("C36").Value = 3.351
$ws1.Range("D36").Value = 48764044.2072
$ws1.Range("E36").Value = 14552087.2
$ws1.Range("F36").Value = 4.165
$ws1.Range("G36").Value = 3.26008629
$ws1.Range("C37").Value = 5.046
$ws1.Range("D37").Value = 40137826.91184
$ws1.Range("E37").Value = 7954385.04
$ws1.Range("F37").Value = 5.081
$ws1.Range("G37").Value = 4.98821585
$ws1.Range("C38").Value = 0.3319
$ws1.Range("D38").Value = 37657550.2389
$ws1.Range("E38").Value = 113460531
$ws1.Range("F38").Value = 5.432
$ws1.Range("G38").Value = 0.32373517
$ws1.Range("A39").Value = 'Raydium'
$ws1.Range("B39").Value = 'RAY'
$ws1.Range("C39").Value = 5.571
$ws1.Range("D39").Value = 38169919.3122
$ws1.Range("E39").Value = 6851538.2
$ws1.Range("F39").Value = 14.16
$ws1.Range("G39").Value = 5.33245647
$ws1.Range("A40").Value = 'Curve DAO'
$ws1.Range("B40").Value = 'CRV'
$ws1.Range("C40").Value = 0.5729
$ws1.Range("D40").Value = 37010394.65160999
$ws1.Range("E40").Value = 64601840.9
$ws1.Range("F40").Value = 9.52
$ws1.Range("G40").Value = 0.55912403
$ws1.Range("C41").Value = 9.819
$ws1.Range("D41").Value = 36517724.18829001
$ws1.Range("E41").Value = 3719087.91
$ws1.Range("F41").Value = 7.558
$ws1.Range("G41").Value = 9.62897315
$ws1.Range("A42").Value = 'Linea Bridged LDO (Linea)'
$ws1.Range("B42").Value = 'LDO'
$ws1.Range("C42").Value = 1.62
$ws1.Range("D42").Value = 35165634.9534
$ws1.Range("E42").Value = 21707182.07
$ws1.Range("F42").Value = 6.649
$ws1.Range("G42").Value = 1.61161134
$ws1.Range("A43").Value = 'Arkham'
$ws1.Range("B43").Value = 'ARKM'
$ws1.Range("C43").Value = 0.707
$ws1.Range("D43").Value = 35051478.2996
$ws1.Range("E43").Value = 49577762.8
$ws1.Range("F43").Value = 6.476
$ws1.Range("G43").Value = 0.68386333
$ws1.Range("C44").Value = 0.2604
$ws1.Range("D44").Value = 31725405.642
$ws1.Range("E44").Value = 121833355
$ws1.Range("F44").Value = 3.952
$ws1.Range("G44").Value = 0.26390724
$ws1.Range("C45").Value = 0.4636
$ws1.Range("D45").Value = 31921531.12412
$ws1.Range("E45").Value = 68855761.7
$ws1.Range("F45").Value = 10.724
$ws1.Range("G45").Value = 0.4457067
$ws1.Range("C46").Value = 0.798
$ws1.Range("D46").Value = 29928895.038
$ws1.Range("E46").Value = 37504881
$ws1.Range("F46").Value = 2.439
$ws1.Range("G46").Value = 0.79236586
$ws1.Range("C47").Value = 0.00009725
$ws1.Range("D47").Value = 29563992.3131655
$ws1.Range("E47").Value = 303999920958
$ws1.Range("F47").Value = 3.689
$ws1.Range("G47").Value = 0.00009692
$ws1.Range("A48").Value = 'Strategic Hub for Innovation in Blockchain'
$ws1.Range("B48").Value = 'SHIB'
$ws1.Range("C48").Value = 0.00001604
$ws1.Range("D48").Value = 28923325.20721144
$ws1.Range("E48").Value = 1803199825886
$ws1.Range("F48").Value = 1.777
$ws1.Range("G48").Value = 0.00001607
$ws1.Range("C49").Value = 0.01597
$ws1.Range("D49").Value = 27633807.05517
$ws1.Range("E49").Value = 1730357361
$ws1.Range("F49").Value = 21.537
$ws1.Range("G49").Value = 0.0166497
$ws1.Range("A50").Value = 'Magic Eden'
$ws1.Range("B50").Value = 'ME'
$ws1.Range("C50").Value = 1.871
$ws1.Range("D50").Value = 28711238.66637
$ws1.Range("E50").Value = 15345397.47
$ws1.Range("F50").Value = 3.142
$ws1.Range("G50").Value = 1.85567692
$ws1.Range("A51").Value = 'Rocket Pool'
$ws1.Range("B51").Value = 'RPL'
$ws1.Range("C51").Value = 11.38
$ws1.Range("D51").Value = 30789149.4408
$ws1.Range("E51").Value = 2705549.16
$ws1.Range("F51").Value = 51.129
$ws1.Range("G51").Value = 10.3128274

# --- Analysis sheet updates ---
$ws2.Range("B1").Value = '2025-02-11 18:05'
$ws2.Range("B3").NumberFormat = "@"
$ws2.Range("B3").Value = '$9,490,989,624.88'
$ws2.Range("B3").Style = "Normal"
$ws2.Range("B4").NumberFormat = "@"
$ws2.Range("B4").Value = '$23,657,869,021,964.33'
$ws2.Range("B4").Style = "Normal"
$ws2.Range("B5").NumberFormat = "@"
$ws2.Range("B5").Value = '$2,043.88'
$ws2.Range("B5").Style = "Normal"
$ws2.Range("B8").Value = 42
$ws2.Range("B9").Value = 8
